# The workbook is a single-sheet "نواقص الأصناف" (item shortages) report.
# Row 7 is the sole data row, and this edit swaps it out for a different
# item/entry while leaving every other cell (headers, dates, footer, etc.)
# untouched:
#
#   C7 (الاسم / item name)             : ATROVENT 250MCG/2ML 20 UNIT DOSE VIAL -> ABIMOL 300MG 5 RECTAL SUPP.
#   H7 (الرصيد الحالي / current stock) : 2:2  -> 33:0
#   L7 (حد الطلب / order limit)        : 1    -> 1   (unchanged)
#   N7 (السعر / price)                 : 286.00   -> 15.00
#   P7 (سعر البيع / sell price)        : 57.2000  -> 15.0000
#   Q7 (عدد التعاملات / tx count)      : 0:4  -> 1:0
#
# All of these cells are stored as literal text in the workbook (not real
# numbers/times), matching their original shared-string representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values: none of these are parsed as numbers by Excel, so a
# direct assignment keeps them as text and leaves the existing cell
# formatting/style completely untouched.
$ws.Range("C7").Value = "ABIMOL 300MG 5 RECTAL SUPP."
$ws.Range("H7").Value = "33:0"
$ws.Range("N7").Value = "15.00"
$ws.Range("Q7").Value = "1:0"

# L7 ("1") and P7 ("15.0000") already carry a numeric cell format
# (#,##0.## style / 0.00), so assigning them directly would make Excel
# reinterpret the text as a genuine number and drop its literal text
# representation. Temporarily force a text format while writing the
# value, then restore the original number format so the cell's style
# (and therefore its font/fill/border/alignment) is left exactly as it
# was.
$fmtL7 = $ws.Range("L7").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $fmtL7

$fmtP7 = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "15.0000"
$ws.Range("P7").NumberFormat = $fmtP7
